$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.250.73'
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("D3").Value = '3.131.56'
$ws.Range("E3").Value = '  +0.57%  '

$ws.Range("E4").Value = '  -0.04%  '

$__style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'581.38"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("E6").Value = '  +0.77%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  -0.28%  '

$ws.Range("E9").Value = '  -0.22%  '

$ws.Range("E10").Value = '  -1.40%  '

$ws.Range("E11").Value = '  -0.09%  '

$__style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.0000250"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = '  +0.14%  '

$__style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'37.57"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = '  +0.86%  '

$__style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.123"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = '  -1.55%  '

$ws.Range("D15").Value = '3.649.57'
$ws.Range("E15").Value = '  +0.59%  '

$ws.Range("D16").Value = '67.156.77'
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("D18").Value = '3.128.93'
$ws.Range("E18").Value = '  +0.49%  '

$__style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'16.40"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = '  +1.34%  '

$__style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'493.81"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = '  +2.13%  '

$ws.Range("E21").Value = '  -0.78%  '

$__style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'7.93"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = '  +5.15%  '

$__style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'84.32"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = '  +0.19%  '

$__style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'13.38"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = '  +1.93%  '

$ws.Range("E25").Value = '  -3.13%  '

$__style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'10.47"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = '  +4.03%  '

$ws.Range("E27").Value = '  +0.00%  '

$__style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'7.96"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("E29").Value = '  -1.77%  '

$__style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'2.70"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = '  +0.21%  '

$__style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'28.79"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("E32").Value = '  -0.46%  '

$ws.Range("D33").Value = '0.0₃0949'
$ws.Range("E33").Value = '  -6.04%  '

$__style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = '  -0.15%  '

$__style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'5.92"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  +0.15%  '

$ws.Range("E36").Value = '  -2.75%  '

$__style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'46.81"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = '  -2.93%  '

$ws.Range("E38").Value = '  -2.62%  '

$__style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'50.19"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("E40").Value = '  -1.48%  '

$ws.Range("E41").Value = '  +1.81%  '

$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("D43").Value = '2.841.36'
$ws.Range("E43").Value = '  +0.12%  '

$__style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'387.28"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = '  +1.40%  '

$__style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'2.62"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = '  -6.64%  '

$ws.Range("E46").Value = '  -1.94%  '

$ws.Range("E47").Value = '  +0.25%  '

$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("E49").Value = '  +0.18%  '

$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("E51").Value = '  -0.27%  '
